# Generate Report for Handoff
# Rolls the Open Localization handback report forward to a newer source
# revision: the tracked markdown file's GUID changes, the generated .xlf
# handoff-package hash changes, and the associated timestamps advance.

$wb = $excel.ActiveWorkbook

$oldGuid = "7312ea04-2ef3-4e98-bca4-af05f5c8dec2"
$newGuid = "8065581c-559f-45a9-a175-93f3cbafd4ca"
$oldHash = "08ad2796595742fccba25ff255346d75b25dd423"
$newHash = "829fbc2bb48ea4a57977d7fccd56083cce33b0ff"

# The external hyperlink target itself does not change -- only the
# human-readable display text (which embeds the old GUID) needs updating.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e67fd7bdfa862b8928743794b93635d8ede03a6/e2e/$oldGuid.md"

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = "2016-08-21 01:02:51"

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkAddress, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 01:02:47"

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkAddress, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 01:02:51"
